# Atualizado por script em 11-11-2023 20:45
# Appends three new match rows (127-129) to the Ekstraklasa 2023-2024 sheet,
# mirroring the formatting of the last existing row (126).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New row data: Indice, pais, liga, temporada, data(serial), casa, golCasa,
# fora, golFora, then 4x (odd, data-odd) pairs, and the match URL.
$rows = @(
    @{
        r = 127
        indice = 126
        data = 45241.625
        casa = "LKS Lodz"; golCasa = 3
        fora = "Piast Gliwice"; golFora = 3
        j = 3.86; k = "05/11/2023 12:42"
        l = 5.26; m = "11/11/2023 14:53"
        n = 3.31; o = "05/11/2023 12:42"
        p = 3.66; q = "11/11/2023 14:53"
        r2 = 2;    s = "05/11/2023 12:42"
        t = 1.75; u = "11/11/2023 14:53"
        v = "https://www.betexplorer.com/football/poland/ekstraklasa/lks-lodz-piast-gliwice/n52Liw0c/"
    },
    @{
        r = 128
        indice = 127
        data = 45241.72916666666
        casa = "Zaglebie"; golCasa = 1
        fora = "Widzew Lodz"; golFora = 1
        j = 1.93; k = "05/11/2023 15:13"
        l = 2.36; m = "11/11/2023 17:21"
        n = 3.61; o = "05/11/2023 15:13"
        p = 3.52; q = "11/11/2023 17:27"
        r2 = 3.74; s = "05/11/2023 15:13"
        t = 3.08; u = "11/11/2023 17:21"
        v = "https://www.betexplorer.com/football/poland/ekstraklasa/zaglebie-widzew-lodz/2ysjxb0A/"
    },
    @{
        r = 129
        indice = 128
        data = 45241.83333333334
        casa = "Korona Kielce"; golCasa = 2
        fora = "Jagiellonia"; golFora = 2
        j = 2.57; k = "04/11/2023 20:12"
        l = 2.68; m = "11/11/2023 19:59"
        n = 3.4;  o = "04/11/2023 20:12"
        p = 3.5;  q = "11/11/2023 19:58"
        r2 = 2.78; s = "04/11/2023 20:12"
        t = 2.67; u = "11/11/2023 19:59"
        v = "https://www.betexplorer.com/football/poland/ekstraklasa/korona-kielce-jagiellonia/tSufyIFG/"
    }
)

foreach ($row in $rows) {
    $rn = $row.r

    # Replicate the exact cell formatting (bold/border style on A, datetime
    # number format on E) from the preceding row before writing values.
    $ws.Range("A126:V126").Copy() | Out-Null
    $ws.Range("A" + $rn + ":V" + $rn).PasteSpecial($xlPasteFormats) | Out-Null

    $ws.Range("A" + $rn).Value = $row.indice
    $ws.Range("B" + $rn).Value = "poland"
    $ws.Range("C" + $rn).Value = "ekstraklasa"
    $ws.Range("D" + $rn).Value = "2023-2024"
    $ws.Range("E" + $rn).Value = $row.data
    $ws.Range("F" + $rn).Value = $row.casa
    $ws.Range("G" + $rn).Value = $row.golCasa
    $ws.Range("H" + $rn).Value = $row.fora
    $ws.Range("I" + $rn).Value = $row.golFora
    $ws.Range("J" + $rn).Value = $row.j
    $ws.Range("K" + $rn).Value = $row.k
    $ws.Range("L" + $rn).Value = $row.l
    $ws.Range("M" + $rn).Value = $row.m
    $ws.Range("N" + $rn).Value = $row.n
    $ws.Range("O" + $rn).Value = $row.o
    $ws.Range("P" + $rn).Value = $row.p
    $ws.Range("Q" + $rn).Value = $row.q
    $ws.Range("R" + $rn).Value = $row.r2
    $ws.Range("S" + $rn).Value = $row.s
    $ws.Range("T" + $rn).Value = $row.t
    $ws.Range("U" + $rn).Value = $row.u
    $ws.Range("V" + $rn).Value = $row.v
}

$excel.CutCopyMode = 0
